$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = "Sunday, Jan 08"
$ws.Range("C34").Value = "8:20 PM"
$ws.Range("D34").Value = "FR6265"
$ws.Range("E34").Value = "Krakow"
$ws.Range("F34").Value = "(KRK)"
$ws.Range("G34").Value = "Ryanair "
$ws.Range("H34").Value = "B738"
$ws.Range("I34").Value = "(SP-RSM)"
$ws.Range("J34").Value = "8:20 PM"
$ws.Range("L34").Value = "0 hours, 0 minutes"

# Row 35
$ws.Range("A35").Value = 34
$ws.Range("B35").Value = "Sunday, Jan 08"
$ws.Range("C35").Value = "10:15 PM"
$ws.Range("D35").Value = "FR7679"
$ws.Range("E35").Value = "Stockholm"
$ws.Range("F35").Value = "(ARN)"
$ws.Range("G35").Value = "Ryanair "
$ws.Range("H35").Value = "B38M"
$ws.Range("I35").Value = "(9H-VUJ)"
$ws.Range("J35").Value = "10:26 PM"
$ws.Range("L35").Value = "0 hours, 11 minutes"
